# Replacing the Kahraman model with a new one.
#
# Column A holds the quarter-hourly timestamps (as Excel serial dates) for
# a single day; the whole series is shifted forward by 4 days (from
# 2025-03-31 to 2025-04-04). Column B holds the predicted solar production
# (MW) for each timestamp; the new model changes the production curve for
# the rows that previously had non-zero morning/daytime output (rows 27-66,
# i.e. 06:30 through 16:15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift every timestamp in column A forward by 4 days -------------------
for ($r = 2; $r -le 97; $r++) {
    $oldDate = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $oldDate + 4
}

# --- Apply the new solar-production curve (column B) -----------------------
$newProduction = @{
    27 = 13
    28 = 35
    29 = 73
    30 = 140
    31 = 187
    32 = 247
    33 = 340
    34 = 402
    35 = 461
    36 = 517
    37 = 582
    38 = 653
    39 = 726
    40 = 777
    41 = 717
    42 = 801
    43 = 870
    44 = 865
    45 = 894
    46 = 899
    47 = 928
    48 = 905
    49 = 899
    50 = 907
    51 = 890
    52 = 947
    53 = 943
    54 = 861
    55 = 859
    56 = 874
    57 = 882
    58 = 891
    59 = 855
    60 = 826
    61 = 757
    62 = 663
    63 = 640
    64 = 599
    65 = 584
    66 = 544
}

foreach ($row in $newProduction.Keys) {
    $ws.Cells.Item($row, 2).Value = $newProduction[$row]
}
